$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" date column (C) from 45207 to 45208 for rows 2-10
for ($row = 2; $row -le 10; $row++) {
    $ws.Cells.Item($row, 3).Value2 = 45208
}
